$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    13011,10599,10599,10599,10599,10599,10208,10208,10176,8959,8959,8959,8959,8959,8959,8959,8959,8959,8888,8225,8225,8225,8225,8225,8225,8225,8225,8027,8027,8027,7909,7909,7909,7909,7909,7909,7909,7909,7909,7639,7639,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
